$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column F with the value "Infra" for every existing data row (1-9)
for ($r = 1; $r -le 9; $r++) {
    $ws.Cells.Item($r, 6).Value = "Infra"
}
